$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Отчет по инвентарю"

# --- Row 1 (top header row) ---
$ws.Range("B1").Value = "Наименование объекта нефинансового учета"
$ws.Range("C1").Value = "Номер(код) объекта учета"
$ws.Range("D1").Value = "Фактическое наличие"
$ws.Range("E1").Value = ""
$ws.Range("F1").Value = ""
$ws.Range("G1").Value = ""
$ws.Range("H1").Value = "По данным бухгалтерского учета"
$ws.Range("I1").Value = ""
$ws.Range("J1").Value = "Результаты инвентаризации"
$ws.Range("K1").Value = ""
$ws.Range("L1").Value = ""
$ws.Range("M1").Value = ""
$ws.Range("N1").Value = "Материально ответственное лицо"

# --- Row 2 ---
$ws.Range("J2").Value = "Отклонение"

# --- Row 3 ---
$ws.Range("D3").Value = "Цена(оценочная стоимость), руб"
$ws.Range("E3").Value = "Количество"
$ws.Range("F3").Value = "Сумма, руб"
$ws.Range("G3").Value = "Статус объекта учета"
$ws.Range("H3").Value = "Количество"
$ws.Range("I3").Value = "Балансовая стоимость, руб"
$ws.Range("J3").Value = "Недосдача"
$ws.Range("L3").Value = "Излишки"

# --- Row 4 ---
$ws.Range("J4").Value = "Количество"
$ws.Range("K4").Value = "Сумма, руб"
$ws.Range("L4").Value = "Количество"
$ws.Range("M4").Value = "Сумма, руб"

# --- Row 5 (was row 2, data moved down + values changed) ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Placeholder"
$ws.Range("C5").Value = "'0"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = "В работе"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = -1
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = "admin"
$ws.Range("O5").Value = "sgdfgd"

# --- Row 6 (new row) ---
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Placeholder"
$ws.Range("C6").Value = "'1"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = "В работе"
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = "admin"
$ws.Range("O6").Value = ""
$ws.Range("O6").Style = "Normal"

# --- Merge cells ---
$ws.Range("A1:A4").Merge()
$ws.Range("B1:B4").Merge()
$ws.Range("C1:C4").Merge()
$ws.Range("D1:G2").Merge()
$ws.Range("H1:I2").Merge()
$ws.Range("J1:M1").Merge()
$ws.Range("N1:N4").Merge()
$ws.Range("O1:O4").Merge()
$ws.Range("J2:M2").Merge()
$ws.Range("D3:D4").Merge()
$ws.Range("E3:E4").Merge()
$ws.Range("F3:F4").Merge()
$ws.Range("G3:G4").Merge()
$ws.Range("H3:H4").Merge()
$ws.Range("I3:I4").Merge()
$ws.Range("J3:K3").Merge()
$ws.Range("L3:M3").Merge()
